$d = $word.ActiveDocument

$find = $d.Content.Find
$find.Execute(
    "Lilith (neutral embarrassed_slightly): This Thursday, okay?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Lilith (neutral embarrassed_slightly): This Thursday, okay?",
    2
)
